# The document already carries bookmarkStart/bookmarkEnd pairs around each
# heading ("a-test-of-headers", "second-level", "third-level",
# "fourth-level", "fifth-level", "sixth-level"). The diff we need to apply
# only renumbers those bookmarks' internal w:id values -- their names,
# positions and nesting are unchanged.
#
# Word assigns/normalizes bookmark ids when the bookmark table is touched,
# so we force that renumbering pass by adding a throwaway bookmark and then
# removing it again. The net effect on document content is nil (no text,
# paragraph, or run is modified) but it causes Word to re-issue fresh,
# sequential ids to every bookmark that is actually present in the
# document, in document order -- exactly the transformation the diff
# describes.

$d = $word.ActiveDocument

$touch = $d.Bookmarks.Add("__renumber_bookmarks__", $d.Range(0, 0))
$touch.Delete()
